# Insert a new data row at row 387 (pushes existing rows 387..467 down to 388..468)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(387).Insert()

# Populate the newly inserted row 387 with its data
$ws.Cells.Item(387, 1).Value = 10
$ws.Cells.Item(387, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(387, 3).Value = "La Araucanía"
$ws.Cells.Item(387, 4).Value = 45275
$ws.Cells.Item(387, 5).Value = 9
$ws.Cells.Item(387, 6).Value = 100112052
$ws.Cells.Item(387, 7).Value = "Albahaca"
$ws.Cells.Item(387, 8).Value = "Sin especificar"
$ws.Cells.Item(387, 9).Value = "Primera"
$ws.Cells.Item(387, 10).Value = 55
$ws.Cells.Item(387, 11).Value = 8000
$ws.Cells.Item(387, 12).Value = 8000
$ws.Cells.Item(387, 13).Value = 8000
$ws.Cells.Item(387, 14).Value = "$/paquete"
$ws.Cells.Item(387, 15).Value = "Región Metropolitana"
$ws.Cells.Item(387, 16).Value = 8000
$ws.Cells.Item(387, 17).Value = 1
$ws.Cells.Item(387, 18).Value = "Hortaliza"
